$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly data point (row for date 45128) is inserted as row 15,
# pushing all existing rows 15-32 down to 16-33.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new record.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 45128
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112013
$ws.Range("G15").Value = "Alcachofa"
$ws.Range("H15").Value = "Española"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = "$/caja 30 unidades"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 483
$ws.Range("Q15").Value = 30
$ws.Range("R15").Value = "Hortaliza"
